$p = $ppt.ActivePresentation

# --- Notes Master: datetimeFigureOut field "7/12/2018" -> "28/9/2018" ---
$nm = $p.NotesMaster
$nm.Shapes.Item(2).TextFrame.TextRange.Text = "28/9/2018"

# --- Slide Layouts (all 11): datetimeFigureOut field "12/7/2018" -> "9/28/2018" ---
$sm = $p.SlideMaster
for ($j = 1; $j -le $sm.CustomLayouts.Count; $j++) {
    $lay = $sm.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.Name.StartsWith("Date Placeholder")) {
            $sh.TextFrame.TextRange.Text = "9/28/2018"
        }
    }
}

# --- Slide 1: "Jest" -> "Karma" ---
$s = $p.Slides.Item(1)
$s.Shapes.Item(25).TextFrame.TextRange.Text = "Karma"
